$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.544.43"
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").Value = "'1.728.34"
$ws.Range("E3").Value = "  -1.18%  "
$ws.Range("D4").Value = "'1.0000"
$ws.Range("E4").Value = "  +0.54%  "
$ws.Range("D5").Value = "'245.24"
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  +0.45%  "
$ws.Range("D7").Value = "'0.4806"
$ws.Range("E7").Value = "  +0.59%  "
$ws.Range("D8").Value = "'0.2672"
$ws.Range("E8").Value = "  -1.62%  "
$ws.Range("D9").Value = "'0.06182"
$ws.Range("E9").Value = "  -1.11%  "
$ws.Range("D10").Value = "'1.731.28"
$ws.Range("E10").Value = "  -0.04%  "
$ws.Range("D11").Value = "'0.07193"
$ws.Range("E11").Value = "  +1.19%  "
$ws.Range("D12").Value = "'15.56"
$ws.Range("E12").Value = "  -2.54%  "
$ws.Range("D13").Value = "'0.6089"
$ws.Range("E13").Value = "  -2.56%  "
$ws.Range("D14").Value = "'4.529"
$ws.Range("E14").Value = "  +0.88%  "
$ws.Range("D15").Value = "'77.26"
$ws.Range("E15").Value = "  +0.20%  "
$ws.Range("D16").Value = "'0.9995"
$ws.Range("E16").Value = "  +0.36%  "
$ws.Range("D17").Value = "'26.556.12"
$ws.Range("E17").Value = "  -0.11%  "
$ws.Range("D18").Value = "'1.001"
$ws.Range("E18").Value = "  +0.59%  "
$ws.Range("D19").Value = "'0.000006953"
$ws.Range("E19").Value = "  +0.72%  "
$ws.Range("D20").Value = "'11.55"
$ws.Range("E20").Value = "  -1.52%  "
$ws.Range("D21").Value = "'1.952.53"
$ws.Range("E21").Value = "  +0.19%  "
$ws.Range("D22").Value = "'4.530"
$ws.Range("E22").Value = "  -1.68%  "
$ws.Range("D23").Value = "'8.803"
$ws.Range("E23").Value = "  -0.78%  "
$ws.Range("D24").Value = "'5.255"
$ws.Range("E24").Value = "  -1.69%  "
$ws.Range("D25").Value = "'137.14"
$ws.Range("E25").Value = "  +0.89%  "
$ws.Range("D26").Value = "'15.39"
$ws.Range("E26").Value = "  -0.31%  "
$ws.Range("D27").Value = "'1.784"
$ws.Range("E27").Value = "  -3.44%  "
$ws.Range("D28").Value = "'1.415"
$ws.Range("E28").Value = "  -0.65%  "
$ws.Range("D29").Value = "'107.69"
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("D30").Value = "'3.978"
$ws.Range("E30").Value = "  -0.81%  "
$ws.Range("D31").Value = "'0.08015"
$ws.Range("E31").Value = "  +1.69%  "
$ws.Range("D32").Value = "'3.691"
$ws.Range("E32").Value = "  -2.09%  "
$ws.Range("D33").Value = "'0.04514"
$ws.Range("E33").Value = "  -1.25%  "
$ws.Range("D34").Value = "'2.614"
$ws.Range("E34").Value = "  -0.33%  "
$ws.Range("D35").Value = "'1.004"
$ws.Range("E35").Value = "  +0.82%  "
$ws.Range("D36").Value = "'0.6279"
$ws.Range("E36").Value = "  -0.92%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "'0.9087"
$ws.Range("E37").Value = "  -5.58%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "'2.068"
$ws.Range("E38").Value = "  +5.63%  "
$ws.Range("D39").Value = "'2.401"
$ws.Range("E39").Value = "  -3.39%  "
$ws.Range("E40").Value = "  +0.41%  "
$ws.Range("D41").Value = "'0.01503"
$ws.Range("E41").Value = "  -0.18%  "
$ws.Range("D42").Value = "'102.23"
$ws.Range("E42").Value = "  -11.33%  "
$ws.Range("D43").Value = "'5.510"
$ws.Range("E43").Value = "  -3.84%  "
$ws.Range("D44").Value = "'0.3890"
$ws.Range("E44").Value = "  -0.80%  "
$ws.Range("D45").Value = "'7.052"
$ws.Range("E45").Value = "  +7.63%  "
$ws.Range("D46").Value = "'0.1181"
$ws.Range("E46").Value = "  -2.26%  "
$ws.Range("D47").Value = "'0.05382"
$ws.Range("E47").Value = "  +1.44%  "
$ws.Range("D48").Value = "'7.868"
$ws.Range("E48").Value = "  -1.69%  "
$ws.Range("D49").Value = "'30.74"
$ws.Range("E49").Value = "  -0.22%  "
$ws.Range("D50").Value = "'1.247"
$ws.Range("E50").Value = "  +0.98%  "
$ws.Range("D51").Value = "'0.3414"
$ws.Range("E51").Value = "  -0.98%  "
